$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(31, 8).Value = 165  # H31
$ws.Cells.Item(31, 9).Value = 165  # I31
$ws.Cells.Item(31, 11).Value = 495  # K31
$ws.Cells.Item(31, 13).Value = -265  # M31
$ws.Cells.Item(86, 8).Value = 8998.166999999999  # H86
$ws.Cells.Item(86, 9).Value = 10495  # I86
$ws.Cells.Item(86, 10).Value = 8249.75  # J86
$ws.Cells.Item(86, 11).Value = 10495  # K86
$ws.Cells.Item(86, 12).Value = 8249.75  # L86
$ws.Cells.Item(86, 13).Value = -9372  # M86
$ws.Cells.Item(86, 14).Value = -10495.75  # N86
$ws.Cells.Item(89, 8).Value = 8998.166999999999  # H89
$ws.Cells.Item(89, 9).Value = 10495  # I89
$ws.Cells.Item(89, 10).Value = 8249.75  # J89
$ws.Cells.Item(89, 11).Value = 52475  # K89
$ws.Cells.Item(89, 12).Value = 41248.75  # L89
$ws.Cells.Item(89, 13).Value = -46859  # M89
$ws.Cells.Item(89, 14).Value = -52480.75  # N89
$ws.Cells.Item(98, 8).Value = 1085.7646  # H98
$ws.Cells.Item(98, 9).Value = 1085.7646  # I98
$ws.Cells.Item(98, 11).Value = 1085.7646  # K98
$ws.Cells.Item(98, 13).Value = 412.2354  # M98
$ws.Cells.Item(116, 8).Value = 17435.334  # H116
$ws.Cells.Item(116, 9).Value = 10000  # I116
$ws.Cells.Item(116, 10).Value = 21153  # J116
$ws.Cells.Item(116, 11).Value = 10000  # K116
$ws.Cells.Item(116, 12).Value = 21153  # L116
$ws.Cells.Item(116, 13).Value = -6558  # M116
$ws.Cells.Item(116, 14).Value = -28037  # N116
$ws.Cells.Item(122, 8).Value = 1085.7646  # H122
$ws.Cells.Item(122, 9).Value = 1085.7646  # I122
$ws.Cells.Item(122, 11).Value = 3257.2938  # K122
$ws.Cells.Item(122, 13).Value = -807.2937999999999  # M122
$ws.Cells.Item(132, 8).Value = 1116.7885  # H132
$ws.Cells.Item(132, 9).Value = 858.6531  # I132
$ws.Cells.Item(132, 11).Value = 2575.9593  # K132
$ws.Cells.Item(132, 13).Value = -45.95929999999998  # M132
$ws.Cells.Item(135, 8).Value = 7144281.5  # H135
$ws.Cells.Item(135, 9).Value = 892.5599999999999  # I135
$ws.Cells.Item(135, 11).Value = 8033.039999999999  # K135
$ws.Cells.Item(135, 13).Value = -5498.039999999999  # M135
$ws.Cells.Item(137, 8).Value = 48783724  # H137
$ws.Cells.Item(137, 9).Value = 35716972  # I137
$ws.Cells.Item(137, 10).Value = 76927496  # J137
$ws.Cells.Item(137, 11).Value = 107150916  # K137
$ws.Cells.Item(137, 12).Value = 230782488  # L137
$ws.Cells.Item(137, 13).Value = -107148366  # M137
$ws.Cells.Item(137, 14).Value = -230787588  # N137
$ws.Cells.Item(138, 8).Value = 7536.2173  # H138
$ws.Cells.Item(138, 9).Value = 3081.7  # I138
$ws.Cells.Item(138, 10).Value = 8773.583000000001  # J138
$ws.Cells.Item(138, 11).Value = 9245.099999999999  # K138
$ws.Cells.Item(138, 12).Value = 26320.749  # L138
$ws.Cells.Item(138, 13).Value = -4105.099999999999  # M138
$ws.Cells.Item(138, 14).Value = -36600.749  # N138

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1245.3334  # H2
$ws.Cells.Item(2, 9).Value = 868.25  # I2
$ws.Cells.Item(2, 10).Value = 1999.5  # J2
$ws.Cells.Item(2, 11).Value = 868.25  # K2
$ws.Cells.Item(2, 12).Value = 1999.5  # L2
$ws.Cells.Item(2, 13).Value = -755.25  # M2
$ws.Cells.Item(2, 14).Value = -2225.5  # N2
$ws.Cells.Item(63, 8).Value = 996.5  # H63
$ws.Cells.Item(63, 9).Value = 999  # I63
$ws.Cells.Item(63, 10).Value = 994  # J63
$ws.Cells.Item(63, 11).Value = 999  # K63
$ws.Cells.Item(63, 12).Value = 994  # L63
$ws.Cells.Item(63, 13).Value = -313  # M63
$ws.Cells.Item(63, 14).Value = -2366  # N63
$ws.Cells.Item(66, 8).Value = 996.5  # H66
$ws.Cells.Item(66, 9).Value = 999  # I66
$ws.Cells.Item(66, 10).Value = 994  # J66
$ws.Cells.Item(66, 11).Value = 4995  # K66
$ws.Cells.Item(66, 12).Value = 4970  # L66
$ws.Cells.Item(66, 13).Value = -1563  # M66
$ws.Cells.Item(66, 14).Value = -11834  # N66
$ws.Cells.Item(116, 8).Value = 1245.3334  # H116
$ws.Cells.Item(116, 9).Value = 868.25  # I116
$ws.Cells.Item(116, 10).Value = 1999.5  # J116
$ws.Cells.Item(116, 11).Value = 868.25  # K116
$ws.Cells.Item(116, 12).Value = 1999.5  # L116
$ws.Cells.Item(116, 13).Value = 1425.75  # M116
$ws.Cells.Item(116, 14).Value = -6587.5  # N116
$ws.Cells.Item(132, 8).Value = 43482700  # H132
$ws.Cells.Item(132, 9).Value = 3905.7222  # I132
$ws.Cells.Item(132, 11).Value = 11717.1666  # K132
$ws.Cells.Item(132, 13).Value = -9187.1666  # M132

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1245.3334  # H3
$ws.Cells.Item(3, 9).Value = 868.25  # I3
$ws.Cells.Item(3, 10).Value = 1999.5  # J3
$ws.Cells.Item(3, 11).Value = 868.25  # K3
$ws.Cells.Item(3, 12).Value = 1999.5  # L3
$ws.Cells.Item(3, 13).Value = -754.25  # M3
$ws.Cells.Item(3, 14).Value = -2227.5  # N3
$ws.Cells.Item(94, 8).Value = 4608.846  # H94
$ws.Cells.Item(94, 9).Value = 1409.5834  # I94
$ws.Cells.Item(94, 11).Value = 1409.5834  # K94
$ws.Cells.Item(94, 13).Value = -958.5834  # M94
$ws.Cells.Item(99, 8).Value = 4122.1113  # H99
$ws.Cells.Item(99, 9).Value = 2388.6316  # I99
$ws.Cells.Item(99, 10).Value = 6059.5293  # J99
$ws.Cells.Item(99, 11).Value = 2388.6316  # K99
$ws.Cells.Item(99, 12).Value = 6059.5293  # L99
$ws.Cells.Item(99, 13).Value = -890.6316000000002  # M99
$ws.Cells.Item(99, 14).Value = -9055.5293  # N99
$ws.Cells.Item(102, 8).Value = 29924.666  # H102
$ws.Cells.Item(102, 9).Value = 9977  # I102
$ws.Cells.Item(102, 10).Value = 69820  # J102
$ws.Cells.Item(102, 11).Value = 9977  # K102
$ws.Cells.Item(102, 12).Value = 69820  # L102
$ws.Cells.Item(102, 13).Value = -6732  # M102
$ws.Cells.Item(102, 14).Value = -76310  # N102
$ws.Cells.Item(134, 8).Value = 1836.25  # H134
$ws.Cells.Item(134, 9).Value = 1836.25  # I134
$ws.Cells.Item(134, 11).Value = 5508.75  # K134
$ws.Cells.Item(134, 13).Value = -2973.75  # M134

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 24394440  # H31
$ws.Cells.Item(31, 9).Value = 3416.5312  # I31
$ws.Cells.Item(31, 10).Value = 111118080  # J31
$ws.Cells.Item(31, 11).Value = 3416.5312  # K31
$ws.Cells.Item(31, 12).Value = 111118080  # L31
$ws.Cells.Item(31, 13).Value = -3121.5312  # M31
$ws.Cells.Item(31, 14).Value = -111118670  # N31
$ws.Cells.Item(34, 8).Value = 24394440  # H34
$ws.Cells.Item(34, 9).Value = 3416.5312  # I34
$ws.Cells.Item(34, 10).Value = 111118080  # J34
$ws.Cells.Item(34, 11).Value = 3416.5312  # K34
$ws.Cells.Item(34, 12).Value = 111118080  # L34
$ws.Cells.Item(34, 13).Value = -3214.5312  # M34
$ws.Cells.Item(34, 14).Value = -111118484  # N34
$ws.Cells.Item(42, 8).Value = 14295  # H42
$ws.Cells.Item(42, 10).Value = 14295  # J42
$ws.Cells.Item(42, 12).Value = 14295  # L42
$ws.Cells.Item(42, 14).Value = -15481  # N42
$ws.Cells.Item(58, 8).Value = 1577  # H58
$ws.Cells.Item(58, 9).Value = 1506.1578  # I58
$ws.Cells.Item(58, 10).Value = 2250  # J58
$ws.Cells.Item(58, 11).Value = 1506.1578  # K58
$ws.Cells.Item(58, 12).Value = 2250  # L58
$ws.Cells.Item(58, 13).Value = -1303.1578  # M58
$ws.Cells.Item(58, 14).Value = -2656  # N58
$ws.Cells.Item(134, 8).Value = 1021.75  # H134
$ws.Cells.Item(134, 9).Value = 1039.8108  # I134
$ws.Cells.Item(134, 11).Value = 3119.4324  # K134
$ws.Cells.Item(134, 13).Value = -584.4323999999997  # M134
$ws.Cells.Item(136, 8).Value = 1577  # H136
$ws.Cells.Item(136, 9).Value = 1506.1578  # I136
$ws.Cells.Item(136, 10).Value = 2250  # J136
$ws.Cells.Item(136, 11).Value = 4518.4734  # K136
$ws.Cells.Item(136, 12).Value = 6750  # L136
$ws.Cells.Item(136, 13).Value = -1968.4734  # M136
$ws.Cells.Item(136, 14).Value = -11850  # N136
$ws.Cells.Item(141, 8).Value = 277720.4  # H141
$ws.Cells.Item(141, 10).Value = 284889.5  # J141
$ws.Cells.Item(141, 12).Value = 284889.5  # L141
$ws.Cells.Item(141, 14).Value = -295249.5  # N141

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(6, 8).Value = 146.35  # H6
$ws.Cells.Item(6, 9).Value = 96.117645  # I6
$ws.Cells.Item(6, 11).Value = 288.352935  # K6
$ws.Cells.Item(6, 13).Value = -175.352935  # M6
$ws.Cells.Item(55, 8).Value = 6252572  # H55
$ws.Cells.Item(55, 10).Value = 3050.3845  # J55
$ws.Cells.Item(55, 12).Value = 9151.1535  # L55
$ws.Cells.Item(55, 14).Value = -9505.1535  # N55
$ws.Cells.Item(98, 8).Value = 3538.4167  # H98
$ws.Cells.Item(98, 10).Value = 3962.4443  # J98
$ws.Cells.Item(98, 12).Value = 11887.3329  # L98
$ws.Cells.Item(98, 14).Value = -14883.3329  # N98
$ws.Cells.Item(128, 8).Value = 116136.25  # H128
$ws.Cells.Item(128, 9).Value = 116136.25  # I128
$ws.Cells.Item(128, 11).Value = 348408.75  # K128
$ws.Cells.Item(128, 13).Value = -343428.75  # M128
$ws.Cells.Item(137, 8).Value = 4513.875  # H137
$ws.Cells.Item(137, 10).Value = 7749.8335  # J137
$ws.Cells.Item(137, 12).Value = 23249.5005  # L137
$ws.Cells.Item(137, 14).Value = -33449.50049999999  # N137

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(11, 8).Value = 30539590  # H11
$ws.Cells.Item(11, 10).Value = 0  # J11
$ws.Cells.Item(11, 12).Value = 0  # L11
$ws.Cells.Item(11, 14).ClearContents()  # delete N11
$ws.Cells.Item(70, 8).Value = 5072.2856  # H70
$ws.Cells.Item(70, 9).Value = 4772.7144  # I70
$ws.Cells.Item(70, 11).Value = 4772.7144  # K70
$ws.Cells.Item(70, 13).Value = -4502.7144  # M70
$ws.Cells.Item(73, 8).Value = 5072.2856  # H73
$ws.Cells.Item(73, 9).Value = 4772.7144  # I73
$ws.Cells.Item(73, 11).Value = 4772.7144  # K73
$ws.Cells.Item(73, 13).Value = -3836.7144  # M73
$ws.Cells.Item(107, 8).Value = 1050.4  # H107
$ws.Cells.Item(107, 9).Value = 1067.1111  # I107
$ws.Cells.Item(107, 10).Value = 900  # J107
$ws.Cells.Item(107, 11).Value = 1067.1111  # K107
$ws.Cells.Item(107, 12).Value = 900  # L107
$ws.Cells.Item(107, 13).Value = 852.8888999999999  # M107
$ws.Cells.Item(107, 14).Value = -4740  # N107
$ws.Cells.Item(132, 8).Value = 2704.1667  # H132
$ws.Cells.Item(132, 9).Value = 2576.5715  # I132
$ws.Cells.Item(132, 10).Value = 3150.75  # J132
$ws.Cells.Item(132, 11).Value = 7729.7145  # K132
$ws.Cells.Item(132, 12).Value = 9452.25  # L132
$ws.Cells.Item(132, 13).Value = -5199.7145  # M132
$ws.Cells.Item(132, 14).Value = -14512.25  # N132

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 5267.364  # H122
$ws.Cells.Item(122, 9).Value = 4657  # I122
$ws.Cells.Item(122, 10).Value = 5999.8  # J122
$ws.Cells.Item(122, 11).Value = 13971  # K122
$ws.Cells.Item(122, 12).Value = 17999.4  # L122
$ws.Cells.Item(122, 13).Value = -11521  # M122
$ws.Cells.Item(122, 14).Value = -22899.4  # N122
$ws.Cells.Item(132, 8).Value = 33901960  # H132
$ws.Cells.Item(132, 9).Value = 3441.8044  # I132
$ws.Cells.Item(132, 10).Value = 153850560  # J132
$ws.Cells.Item(132, 11).Value = 10325.4132  # K132
$ws.Cells.Item(132, 12).Value = 461551680  # L132
$ws.Cells.Item(132, 13).Value = -7795.413199999999  # M132
$ws.Cells.Item(132, 14).Value = -461556740  # N132

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(20, 8).Value = 6748.75  # H20
$ws.Cells.Item(20, 9).Value = 5665  # I20
$ws.Cells.Item(20, 10).Value = 10000  # J20
$ws.Cells.Item(20, 11).Value = 5665  # K20
$ws.Cells.Item(20, 12).Value = 10000  # L20
$ws.Cells.Item(20, 13).Value = -5425  # M20
$ws.Cells.Item(20, 14).Value = -10480  # N20
$ws.Cells.Item(132, 8).Value = 4981.783  # H132
$ws.Cells.Item(132, 9).Value = 4895.22  # I132
$ws.Cells.Item(132, 11).Value = 14685.66  # K132
$ws.Cells.Item(132, 13).Value = -12155.66  # M132
$ws.Cells.Item(136, 8).Value = 1703.381  # H136
$ws.Cells.Item(136, 9).Value = 1068.5333  # I136
$ws.Cells.Item(136, 10).Value = 3290.5  # J136
$ws.Cells.Item(136, 11).Value = 3205.5999  # K136
$ws.Cells.Item(136, 12).Value = 9871.5  # L136
$ws.Cells.Item(136, 13).Value = -655.5999000000002  # M136
$ws.Cells.Item(136, 14).Value = -14971.5  # N136

